$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: replace numeric 4s with shared-string "4L"
$ws.Range("C19:O19").Value = "4L"
$ws.Range("Q19:T19").Value = "4L"
$ws.Range("V19:AO19").Value = "4L"
$ws.Range("AQ19:AQ19").Value = "4L"
$ws.Range("AS19:BM19").Value = "4L"
$ws.Range("BO19:BR19").Value = "4L"
$ws.Range("BT19:CH19").Value = "4L"

# Row 21: replace numeric 4s with shared-string "4R"
$ws.Range("C21:T21").Value = "4R"
$ws.Range("V21:Y21").Value = "4R"
$ws.Range("AA21:AO21").Value = "4R"
$ws.Range("AQ21:AQ21").Value = "4R"
$ws.Range("AS21:BH21").Value = "4R"
$ws.Range("BJ21:BM21").Value = "4R"
$ws.Range("BO21:CH21").Value = "4R"

# Update the view state to match the author's final selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("CK31").Select()
